# ---------------------------------------------------------------------------
# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" worksheet (fund-holding detail) positioned right
#    before the "总计" (totals) sheet.
# 2. Update the "总计" sheet with a new leading row for 2022-Q1, shifting the
#    existing quarters down by one row.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# A template sheet that already has the exact styles we need to reuse
# (header style, index-column style, page margins, outline settings, ...).
$template = $wb.Worksheets.Item("2021-Q4")

# Capture the "总计" sheet's existing data rows before rebuilding it, so the
# script stays correct even if those values were ever to change upstream.
$totalSheetOld = $wb.Worksheets.Item("总计")
$oldRowCount = $totalSheetOld.Cells.Item($totalSheetOld.Rows.Count, 2).End(-4162).Row
$existingQuarters = New-Object System.Collections.ArrayList
for ($r = 2; $r -le $oldRowCount; $r++) {
    $q = $totalSheetOld.Cells.Item($r, 2).Value()
    $cnt = $totalSheetOld.Cells.Item($r, 3).Value()
    $val = $totalSheetOld.Cells.Item($r, 4).Value()
    [void]$existingQuarters.Add(@($q, $cnt, $val))
}

# Delete the old "总计" sheet. Recreating it after the new sheet keeps the
# internal sheetId / relationship-id allocation in the same order a
# from-scratch export would use (2022-Q1 = 6, 总计 = 7).
$totalSheetOld.Delete()

# ---------------------------------------------------------------------------
# Step 1: new "2022-Q1" worksheet
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "2022-Q1"

# sheet-level cosmetics (outline + page margins) matching sibling sheets
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# header row + index-column styles, copied from the template sheet
$template.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2:A8").Copy()
$newSheet.Range("A2:A8").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$fundRows = @(
    @("009562", "工银瑞信中国机会全球配置股票(QDII)美元", "6.65", "92.85", "3.34", "0.2221", 2),
    @("486001", "工银瑞信中国机会全球配置股票(QDII)",     "6.65", "92.85", "3.34", "0.2221", 2),
    @("009563", "工银瑞信中国机会全球配置股票(QDII)港币", "6.65", "92.85", "3.34", "0.2221", 2),
    @("486002", "工银全球精选股票(QDII)",                 "4.23", "94.60", "4.67", "0.1975", 1),
    @("012751", "建信纳斯达克100指数（QDII）A 美元现汇",   "0.34", "88.02", "10.75", "0.0366", 2),
    @("012752", "建信纳斯达克100指数（QDII）C 人民币",     "0.34", "88.02", "10.75", "0.0366", 2),
    @("012753", "建信纳斯达克100指数（QDII）C 美元现汇",   "0.34", "88.02", "10.75", "0.0366", 2)
)

for ($i = 0; $i -lt $fundRows.Count; $i++) {
    $r = 2 + $i
    $row = $fundRows[$i]

    $newSheet.Cells.Item($r, 1).Value = $i

    $newSheet.Cells.Item($r, 2).Value = "'" + $row[0]
    $newSheet.Cells.Item($r, 2).Style = "Normal"

    $newSheet.Cells.Item($r, 3).Value = $row[1]

    $newSheet.Cells.Item($r, 4).Value = "'" + $row[2]
    $newSheet.Cells.Item($r, 4).Style = "Normal"

    $newSheet.Cells.Item($r, 5).Value = "'" + $row[3]
    $newSheet.Cells.Item($r, 5).Style = "Normal"

    $newSheet.Cells.Item($r, 6).Value = "'" + $row[4]
    $newSheet.Cells.Item($r, 6).Style = "Normal"

    $newSheet.Cells.Item($r, 7).Value = "'" + $row[5]
    $newSheet.Cells.Item($r, 7).Style = "Normal"

    $newSheet.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------------
# Step 2: rebuild "总计" right after "2022-Q1", with the new leading row
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Add($null, $newSheet)
$totalSheet.Name = "总计"

$totalSheet.Outline.SummaryRow = 1
$totalSheet.Outline.SummaryColumn = 1
$totalSheet.PageSetup.LeftMargin = 54
$totalSheet.PageSetup.RightMargin = 54
$totalSheet.PageSetup.TopMargin = 72
$totalSheet.PageSetup.BottomMargin = 72
$totalSheet.PageSetup.HeaderMargin = 36
$totalSheet.PageSetup.FooterMargin = 36

$template.Range("B1:D1").Copy()
$totalSheet.Range("B1:D1").PasteSpecial(-4122)

$totalRowCount = $existingQuarters.Count + 1
$template.Range("A2:A8").Copy()
$totalSheet.Range("A2:A" + (1 + $totalRowCount)).PasteSpecial(-4122)

$totalSheet.Range("B1").Value = "日期"
$totalSheet.Range("C1").Value = "持有数量(只)"
$totalSheet.Range("D1").Value = "持有市值(亿元)"

# new leading row for the just-added quarter
$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 7
$totalSheet.Cells.Item(2, 4).Value = 0.97

# the previously-existing quarters, shifted down by one row
for ($i = 0; $i -lt $existingQuarters.Count; $i++) {
    $r = 3 + $i
    $trio = $existingQuarters[$i]
    $totalSheet.Cells.Item($r, 1).Value = $i + 1
    $totalSheet.Cells.Item($r, 2).Value = $trio[0]
    $totalSheet.Cells.Item($r, 3).Value = $trio[1]
    $totalSheet.Cells.Item($r, 4).Value = $trio[2]
}

Write-Output "2022-Q1 sheet added; 总计 sheet updated"
